# Apply the "LinuxForHealth" rebrand / republish edits to the workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-health-insurance-oversight-system-plan"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Fixed Value for Extension.url row now points at the new URL as well
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-health-insurance-oversight-system-plan"

# Constraint(s) cell for the top-level Extension row is cleared out
$elements.Range("AI2").Value = ""
